# parkrun-times.xlsx — replace the elapsed-time ([h]:mm:ss) values in column D
# with plain "minutes.seconds" numbers, reset their number format back to
# General, and restore the sheet's view state (zoom / selection).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Select()

# New D2:D50 values (row -> value), taken from the target workbook.
$newTimes = @{
    2  = 35.5
    3  = 34.08
    4  = 32.32
    5  = 33.38
    6  = 31.4
    7  = 31.12
    8  = 31.06
    9  = 31.36
    10 = 32.22
    11 = 32.04
    12 = 32.22
    13 = 31.48
    14 = 32.31
    15 = 32.08
    16 = 32.19
    17 = 31.16
    18 = 30.57
    19 = 30.57
    20 = 30.47
    21 = 29.44
    22 = 31.38
    23 = 31.36
    24 = 31.15
    25 = 30.58
    26 = 31.34
    27 = 31.29
    28 = 31.55
    29 = 31.38
    30 = 33.02
    31 = 30.04
    32 = 31.54
    33 = 32.18
    34 = 30.35
    35 = 29.35
    36 = 30.2
    37 = 31.15
    38 = 30.54
    39 = 30.17
    40 = 32.38
    41 = 30.02
    42 = 29.14
    43 = 29.21
    44 = 29.45
    45 = 30.19
    46 = 28.53
    47 = 28.28
    48 = 27.25
    49 = 28.14
    50 = 27.59
}

foreach ($row in 2..50) {
    $cell = $ws.Cells.Item($row, 4)
    $cell.Value = $newTimes[$row]
    # The old [h]:mm:ss elapsed-time format no longer applies to these plain
    # numbers - put the cells back on the workbook's default (General) style.
    $cell.Style = "Normal"
}

# Restore the sheet view: smaller zoom, scrolled down a bit, new selection.
$excel.ActiveWindow.Zoom = 75
$excel.ActiveWindow.ScrollRow = 19
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("F49").Select()
